{"js": "// Replace the 25 three-digit-by-one-digit multiplication expressions in the\n// table cells with their new values, per the commit's regenerated answers.\nconst replacements = [\n  [\"285\\u00D75=1425\", \"639\\u00D78=5112\"],\n  [\"864\\u00D79=7776\", \"370\\u00D74=1480\"],\n  [\"615\\u00D78=4920\", \"748\\u00D72=1496\"],\n  [\"211\\u00D74=844\", \"583\\u00D72=1166\"],\n  [\"855\\u00D75=4275\", \"562\\u00D75=2810\"],\n  [\"462\\u00D74=1848\", \"506\\u00D76=3036\"],\n  [\"415\\u00D76=2490\", \"545\\u00D75=2725\"],\n  [\"983\\u00D75=4915\", \"332\\u00D72=664\"],\n  [\"213\\u00D79=1917\", \"794\\u00D72=1588\"],\n  [\"428\\u00D78=3424\", \"605\\u00D75=3025\"],\n  [\"567\\u00D78=4536\", \"843\\u00D76=5058\"],\n  [\"718\\u00D76=4308\", \"359\\u00D74=1436\"],\n  [\"194\\u00D76=1164\", \"202\\u00D76=1212\"],\n  [\"401\\u00D78=3208\", \"866\\u00D72=1732\"],\n  [\"780\\u00D76=4680\", \"915\\u00D77=6405\"],\n  [\"841\\u00D77=5887\", \"201\\u00D79=1809\"],\n  [\"554\\u00D76=3324\", \"489\\u00D77=3423\"],\n  [\"998\\u00D77=6986\", \"909\\u00D79=8181\"],\n  [\"717\\u00D76=4302\", \"606\\u00D74=2424\"],\n  [\"221\\u00D76=1326\", \"102\\u00D79=918\"],\n  [\"135\\u00D76=810\", \"655\\u00D74=2620\"],\n  [\"823\\u00D75=4115\", \"474\\u00D75=2370\"],\n  [\"820\\u00D76=4920\", \"554\\u00D72=1108\"],\n  [\"442\\u00D72=884\", \"216\\u00D72=432\"],\n  [\"775\\u00D75=3875\", \"487\\u00D76=2922\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Build a quick lookup of old -> new text.\nconst map = new Map(replacements);\n\nfor (const paragraph of paragraphs.items) {\n  // Paragraph.text in Office.js excludes the trailing paragraph mark, so this\n  // compares cleanly against the OOXML w:t contents (each cell has a single\n  // run holding the whole expression).\n  const current = paragraph.text;\n  if (map.has(current)) {\n    const newText = map.get(current);\n    paragraph.getRange().insertText(newText, Word.InsertLocation.replace);\n    map.delete(current);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 three-digit-by-one-digit multiplication expressions in the\n# table cells with their new values, per the commit's regenerated answers.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"285\u00d75=1425\", \"639\u00d78=5112\"),\n    @(\"864\u00d79=7776\", \"370\u00d74=1480\"),\n    @(\"615\u00d78=4920\", \"748\u00d72=1496\"),\n    @(\"211\u00d74=844\", \"583\u00d72=1166\"),\n    @(\"855\u00d75=4275\", \"562\u00d75=2810\"),\n    @(\"462\u00d74=1848\", \"506\u00d76=3036\"),\n    @(\"415\u00d76=2490\", \"545\u00d75=2725\"),\n    @(\"983\u00d75=4915\", \"332\u00d72=664\"),\n    @(\"213\u00d79=1917\", \"794\u00d72=1588\"),\n    @(\"428\u00d78=3424\", \"605\u00d75=3025\"),\n    @(\"567\u00d78=4536\", \"843\u00d76=5058\"),\n    @(\"718\u00d76=4308\", \"359\u00d74=1436\"),\n    @(\"194\u00d76=1164\", \"202\u00d76=1212\"),\n    @(\"401\u00d78=3208\", \"866\u00d72=1732\"),\n    @(\"780\u00d76=4680\", \"915\u00d77=6405\"),\n    @(\"841\u00d77=5887\", \"201\u00d79=1809\"),\n    @(\"554\u00d76=3324\", \"489\u00d77=3423\"),\n    @(\"998\u00d77=6986\", \"909\u00d79=8181\"),\n    @(\"717\u00d76=4302\", \"606\u00d74=2424\"),\n    @(\"221\u00d76=1326\", \"102\u00d79=918\"),\n    @(\"135\u00d76=810\", \"655\u00d74=2620\"),\n    @(\"823\u00d75=4115\", \"474\u00d75=2370\"),\n    @(\"820\u00d76=4920\", \"554\u00d72=1108\"),\n    @(\"442\u00d72=884\", \"216\u00d72=432\"),\n    @(\"775\u00d75=3875\", \"487\u00d76=2922\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
